# TrialsSetup.xlsx update (2026-02-06 12:00)
# The "REFRACT" trial has completed / dropped out of the portfolio, so its
# row is removed from the refreshed query results. The remaining trials'
# "Days remaining" counters have also ticked down by one day since the
# last refresh (REJOICE: 14 -> 13, REMASTER: 34 -> 33).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the REFRACT row (row 2) entirely - this shifts every following
# row up by one and shrinks the table/used-range automatically.
$ws.Rows(2).Delete()

# Keep the hidden ExternalData_1 defined name (used by the query table)
# in sync with the now-smaller range.
$wb.Names("ExternalData_1").RefersTo = "=Sheet1!`$A`$1:`$C`$10"

# Refreshed "Days remaining" values for REJOICE and REMASTER.
$ws.Range("B7").Value = 13
$ws.Range("B9").Value = 33
